# Kehoach.xlsx update:
# - Clear the "x" markers in D16 and C17 (task progress moved)
# - Mark G17/H17 (actual start/end dates) as completed on 19/10/2018
# - Update the active selection to H17

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 ("Thuc hien ham them xoa sua lien ket bang bam"): remove the "x" from D16
$ws.Range("D16").ClearContents()

# Row 17 ("Thuc hien chuc nang tra tu lien ket voi bang bam"): remove the "x" from C17
$ws.Range("C17").ClearContents()

# Row 17: set the actual start/end dates to 19/10/2018
$ws.Range("G17").Value = "19/10/2018"
$ws.Range("H17").Value = "19/10/2018"

# Update the sheet's active selection/cell to H17
$ws.Range("H17").Select()
